# edit.ps1
# Applies the "Updated cryptos list" data refresh described by the commit diff.
# Every Price (D) / Volume(1h) (E) / Coin (B) / Link (C) cell listed below is
# rewritten to its new value. Because several Price cells look like plain
# numbers (e.g. "176.50", "0.517") but must stay literal TEXT (matching the
# workbook's original inlineStr text cells, with exact formatting such as
# trailing zeros preserved), each write briefly flips the cell to the "Text"
# number format, assigns the value, then restores the "Normal" style so the
# cell ends up back at the default (unstyled) format - same as before the
# edit - with the content stored as text instead of being auto-coerced to a
# number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "69.594.89"
Set-TextValue "E2" "  +2.60%  "
Set-TextValue "D3" "2.520.92"
Set-TextValue "E4" "  +0.12%  "
Set-TextValue "D5" "598.85"
Set-TextValue "E5" "  +1.82%  "
Set-TextValue "D6" "176.50"
Set-TextValue "E6" "  -0.35%  "
Set-TextValue "E7" "  -0.02%  "
Set-TextValue "D8" "0.517"
Set-TextValue "E8" "  +0.58%  "
Set-TextValue "D9" "2.520.51"
Set-TextValue "E9" "  +0.88%  "
Set-TextValue "E10" "  +12.40%  "
Set-TextValue "E11" "  -0.28%  "
Set-TextValue "D12" "0.343"
Set-TextValue "E12" "  +1.13%  "
Set-TextValue "D13" "5.02"
Set-TextValue "E13" "  +1.48%  "
Set-TextValue "B14" "ShibaInu"
Set-TextValue "C14" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D14" "0.0000182"
Set-TextValue "E14" "  +6.10%  "
Set-TextValue "B15" "WrappedliquidstakedEther2.0"
Set-TextValue "C15" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue "D15" "2.983.23"
Set-TextValue "E15" "  +1.03%  "
Set-TextValue "D16" "26.02"
Set-TextValue "E16" "  +1.49%  "
Set-TextValue "D17" "69.528.20"
Set-TextValue "D18" "2.498.65"
Set-TextValue "E18" "  +0.86%  "
Set-TextValue "D19" "7.64"
Set-TextValue "E19" "  +2.02%  "
Set-TextValue "D20" "363.94"
Set-TextValue "E20" "  +3.30%  "
Set-TextValue "D21" "11.03"
Set-TextValue "E21" "  +0.51%  "
Set-TextValue "D22" "4.05"
Set-TextValue "E22" "  -1.54%  "
Set-TextValue "E23" "  -0.13%  "
Set-TextValue "D24" "70.45"
Set-TextValue "E24" "  -0.53%  "
Set-TextValue "D25" "4.24"
Set-TextValue "E25" "  -1.23%  "
Set-TextValue "B26" "Aptos"
Set-TextValue "C26" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D26" "9.14"
Set-TextValue "E26" "  +0.22%  "
Set-TextValue "B27" "SuiNetwork"
Set-TextValue "C27" "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue "D27" "1.71"
Set-TextValue "E27" "  -1.94%  "
Set-TextValue "D28" "2.647.63"
Set-TextValue "E28" "  +2.33%  "
Set-TextValue "D29" "0.998"
Set-TextValue "E29" "  -0.04%  "
Set-TextValue "D30" "0.0₃0901"
Set-TextValue "E30" "  -0.84%  "
Set-TextValue "D31" "510.46"
Set-TextValue "E31" "  +1.26%  "
Set-TextValue "D32" "7.72"
Set-TextValue "E32" "  -1.45%  "
Set-TextValue "E33" "  -1.76%  "
Set-TextValue "D34" "1.79"
Set-TextValue "E34" "  +0.79%  "
Set-TextValue "E35" "  +0.13%  "
Set-TextValue "E36" "  -1.30%  "
Set-TextValue "D37" "161.36"
Set-TextValue "E37" "  -1.80%  "
Set-TextValue "D38" "18.74"
Set-TextValue "E38" "  +2.08%  "
Set-TextValue "E39" "  +1.38%  "
Set-TextValue "E40" "  -1.41%  "
Set-TextValue "E41" "  +0.03%  "
Set-TextValue "D42" "1.73"
Set-TextValue "E42" "  -0.57%  "
Set-TextValue "D43" "4.80"
Set-TextValue "E43" "  -1.37%  "
Set-TextValue "E44" "  -2.46%  "
Set-TextValue "D45" "2.39"
Set-TextValue "E45" "  -3.04%  "
Set-TextValue "D46" "38.79"
Set-TextValue "E46" "  -0.47%  "
Set-TextValue "D47" "151.26"
Set-TextValue "E47" "  +4.31%  "
Set-TextValue "D48" "3.59"
Set-TextValue "E48" "  +1.49%  "
Set-TextValue "D49" "0.517"
Set-TextValue "E49" "  -0.14%  "
Set-TextValue "D50" "0.0739"
Set-TextValue "E50" "  -0.63%  "
Set-TextValue "E51" "  -1.68%  "

Write-Host "Applied cryptos list update."
